# Update Name of Algo
# Applies refreshed KNN-imputed values to columns D and E (terrestrial_mammals,
# combination_2_ABCDE/DE/20/seed2) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.682
$ws.Range("D3").Value = -8.453999999999999
$ws.Range("D14").Value = -8.005000000000001
$ws.Range("D16").Value = -8.128
$ws.Range("E18").Value = 16.371
$ws.Range("D21").Value = -8.440999999999999
$ws.Range("D23").Value = -7.874
$ws.Range("E24").Value = 16.763
$ws.Range("D25").Value = -7.842999999999999
$ws.Range("E25").Value = 17.239
$ws.Range("D26").Value = -7.925999999999999
$ws.Range("E27").Value = 16.844
$ws.Range("D29").Value = -7.404999999999999
$ws.Range("E30").Value = 17.074
$ws.Range("E31").Value = 17.118
$ws.Range("E39").Value = 16.573
$ws.Range("D40").Value = -8.273999999999999
$ws.Range("E42").Value = 16.663
$ws.Range("E48").Value = 17.179
$ws.Range("E51").Value = 16.443
$ws.Range("E52").Value = 16.543
$ws.Range("D53").Value = -7.784999999999999
$ws.Range("E55").Value = 16.508
$ws.Range("E56").Value = 16.214
$ws.Range("D57").Value = -7.904000000000001
$ws.Range("E57").Value = 16.512
$ws.Range("D59").Value = -8.098000000000001
$ws.Range("E60").Value = 16.592
$ws.Range("D65").Value = -7.824
$ws.Range("D69").Value = -7.597
$ws.Range("E73").Value = 16.572
$ws.Range("E74").Value = 16.63
$ws.Range("D79").Value = -7.885
$ws.Range("D83").Value = -8.347
$ws.Range("E89").Value = 17.57
$ws.Range("E90").Value = 16.73
$ws.Range("D91").Value = -7.187
$ws.Range("E92").Value = 17.253
$ws.Range("D93").Value = -7.658000000000001
$ws.Range("D100").Value = -7.925999999999999

$wb.Save()
